# Update the code-extractor result table (rows 2-9) to reflect the
# corrected test data: rename the Java controller test id so it no
# longer collides between the "Test" and "Tst" prefixed variants, and
# shift every other row down to keep the table internally consistent.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'TestHelloController.java_Test_AddMethod'
$ws.Range("B2").Value = @'
public void Test_AddMethod() {  
            BasicMaths bm = new BasicMaths();  
            double res = bm.Add(1, 10);  
            Assert.AreEqual(res, 5);  
        } 

'@
$ws.Rows.Item(2).AutoFit()

$ws.Range("A3").Value = 'TestHelloController.java_Test_AddMethod'
$ws.Range("B3").Value = @'
public void Test_AddMethod() {  
            BasicMaths bm = new BasicMaths();  
            double res = bm.Add(10, 10);  
            Assert.AreEqual(res, 2);  
        } 

'@
$ws.Rows.Item(3).AutoFit()

$ws.Range("A4").Value = 'TstHelloController.java_Test_AddMethod'
$ws.Range("B4").Value = @'
public void Test_AddMethod() {  
            BasicMaths bm = new BasicMaths();  
            double res = bm.Add(1, 10);  
            Assert.AreEqual(res, 5);  
        } 

'@
$ws.Rows.Item(4).AutoFit()

$ws.Range("A5").Value = 'TstHelloController.java_Test_AddMethod'
$ws.Range("B5").Value = @'
public void Test_AddMethod() {  
            BasicMaths bm = new BasicMaths();  
            double res = bm.Add(10, 10);  
            Assert.AreEqual(res, 2);  
        } 

'@
$ws.Rows.Item(5).AutoFit()

$ws.Range("A6").Value = 'test.cs_Test_AddMethod'
$ws.Range("B6").Value = @'
    public void Test_AddMethod() {  
            BasicMaths bm = new BasicMaths();  
            double res = bm.Add(10, 10);  
            Assert.AreEqual(res, 20);  
        }  

'@
$ws.Rows.Item(6).AutoFit()

$ws.Range("A7").Value = 'test.cs_Test_DivideMethod'
$ws.Range("B7").Value = @'
    public void Test_DivideMethod() {  
            BasicMaths bm = new BasicMaths();  
            double res = bm.divide(10, 5);  
            Assert.AreEqual(res, 2);  
        }  

'@
$ws.Rows.Item(7).AutoFit()

$ws.Range("A8").Value = 'test.cs_Test_MultiplyMethod'
$ws.Range("B8").Value = @'
    public void Test_MultiplyMethod() {  
        BasicMaths bm = new BasicMaths();  
        double res = bm.Multiply(10, 10);  
        Assert.AreEqual(res, 100);  
    }  

'@
$ws.Rows.Item(8).AutoFit()

$ws.Range("A9").Value = 'test.cs_Test_SubstractMethod'
$ws.Range("B9").Value = @'
    public void Test_SubstractMethod() {  
            BasicMaths bm = new BasicMaths();  
            double res = bm.Substract(10, 10);  
            Assert.AreEqual(res, 0);  
        }  

'@
$ws.Rows.Item(9).AutoFit()

$ws.Columns.Item(1).ColumnWidth = 77.5

[void]$ws.Range("A2:A9").Select()

